$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formats from column E into the new column D for all populated row ranges
# (done in separate contiguous blocks so rows with no D/E cell at all, like 37 and 79, are left untouched)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the newest reporting period values (FY2018 / period ending 2018-12-31)
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1038300
$ws.Range("D9").Value = 743600
$ws.Range("D10").Value = 294600
$ws.Range("D12").Value = 98700
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 26200
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 976900
$ws.Range("D18").Value = 61300
$ws.Range("D20").Value = 1700
$ws.Range("D21").Value = 113700
$ws.Range("D22").Value = 4900
$ws.Range("D23").Value = 58100
$ws.Range("D24").Value = 16200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 41900
$ws.Range("D27").Value = 41900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -1700
$ws.Range("D33").Value = 41900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 41900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 39600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 166900
$ws.Range("D44").Value = 112500
$ws.Range("D45").Value = 124100
$ws.Range("D46").Value = 443100
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 171400
$ws.Range("D49").Value = 111700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 76900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 803000
$ws.Range("D57").Value = 93100
$ws.Range("D58").Value = 3400
$ws.Range("D59").Value = 78900
$ws.Range("D60").Value = 175400
$ws.Range("D61").Value = 136500
$ws.Range("D62").Value = 11500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 323300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 364000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 479700
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 41900
$ws.Range("D83").Value = 50600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 118400
$ws.Range("D91").Value = -41500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -40800
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -139300
$ws.Range("D101").Value = -2000
$ws.Range("D102").Value = -63600

# A handful of prior-period (now column E, formerly D) figures were also revised
$ws.Range("E9").Value = 674800
$ws.Range("E10").Value = 310900
$ws.Range("E17").Value = 888600
$ws.Range("E18").Value = 97100
$ws.Range("E20").Value = -23000
$ws.Range("E26").Value = 55300
$ws.Range("E27").Value = 55300
$ws.Range("E29").Value = -20100
$ws.Range("E32").Value = 23000
